$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 148.5
$ws.Range("I9").Value = 83.57143000000001
$ws.Range("J9").Value = 300
$ws.Range("K9").Value = 83.57143000000001
$ws.Range("L9").Value = 300
$ws.Range("M9").Value = 85.42856999999999
$ws.Range("N9").Value = -638

$ws.Range("H76").Value = 6176021.5
$ws.Range("I76").Value = 3366.5
$ws.Range("J76").Value = 7939637
$ws.Range("K76").Value = 3366.5
$ws.Range("L76").Value = 7939637
$ws.Range("M76").Value = -3051.5
$ws.Range("N76").Value = -7940267

$ws.Range("H79").Value = 6176021.5
$ws.Range("I79").Value = 3366.5
$ws.Range("J79").Value = 7939637
$ws.Range("K79").Value = 3366.5
$ws.Range("L79").Value = 7939637
$ws.Range("M79").Value = -2274.5
$ws.Range("N79").Value = -7941821

$ws.Range("H86").Value = 7943.15
$ws.Range("I86").Value = 2840.2
$ws.Range("K86").Value = 2840.2
$ws.Range("M86").Value = -1717.2

$ws.Range("H89").Value = 7943.15
$ws.Range("I89").Value = 2840.2
$ws.Range("K89").Value = 14201
$ws.Range("M89").Value = -8585

$ws.Range("H92").Value = 76923720
$ws.Range("I92").Value = 111111990
$ws.Range("J92").Value = 116.25
$ws.Range("K92").Value = 111111990
$ws.Range("L92").Value = 116.25
$ws.Range("M92").Value = -111110742
$ws.Range("N92").Value = -2612.25

$ws.Range("H106").Value = 10418660
$ws.Range("I106").Value = 14494203
$ws.Range("K106").Value = 14494203
$ws.Range("M106").Value = -14493572

$ws.Range("H129").Value = 1000522.5
$ws.Range("J129").Value = 2000705.6
$ws.Range("L129").Value = 6002116.800000001
$ws.Range("N129").Value = -6012116.800000001

$ws.Range("H137").Value = 2144.0688
$ws.Range("I137").Value = 2069.7273
$ws.Range("K137").Value = 6209.1819
$ws.Range("M137").Value = -3659.1819

$ws.Range("H138").Value = 2174.234
$ws.Range("I138").Value = 1808.4
$ws.Range("J138").Value = 2302.5964
$ws.Range("K138").Value = 5425.200000000001
$ws.Range("L138").Value = 6907.789199999999
$ws.Range("M138").Value = -285.2000000000007
$ws.Range("N138").Value = -17187.7892

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 31249998
$ws.Range("I63").Value = 31249996
$ws.Range("K63").Value = 31249996
$ws.Range("M63").Value = -31249310

$ws.Range("H66").Value = 31249998
$ws.Range("I66").Value = 31249996
$ws.Range("K66").Value = 156249980
$ws.Range("M66").Value = -156246548

$ws.Range("H122").Value = 1942.1578
$ws.Range("I122").Value = 1594.2354
$ws.Range("K122").Value = 4782.706200000001
$ws.Range("M122").Value = -2332.706200000001

$ws.Range("H132").Value = 14702.462
$ws.Range("I132").Value = 1849.4286
$ws.Range("J132").Value = 47419.273
$ws.Range("K132").Value = 5548.2858
$ws.Range("L132").Value = 142257.819
$ws.Range("M132").Value = -3018.2858
$ws.Range("N132").Value = -147317.819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 10000
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 10000
$ws.Range("K12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").Value = 10000
$ws.Range("N12").Value = -10336

$ws.Range("H94").Value = 2653.842
$ws.Range("I94").Value = 2527.5334
$ws.Range("K94").Value = 2527.5334
$ws.Range("M94").Value = -2076.5334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1004
$ws.Range("I13").Value = 1004
$ws.Range("K13").Value = 1004
$ws.Range("M13").Value = -865

$ws.Range("H31").Value = 14229.967
$ws.Range("I31").Value = 23148.334
$ws.Range("K31").Value = 23148.334
$ws.Range("M31").Value = -22853.334

$ws.Range("H34").Value = 14229.967
$ws.Range("I34").Value = 23148.334
$ws.Range("K34").Value = 23148.334
$ws.Range("M34").Value = -22946.334

$ws.Range("H107").Value = 561.5
$ws.Range("I107").Value = 222
$ws.Range("J107").Value = 1014.1667
$ws.Range("K107").Value = 222
$ws.Range("L107").Value = 1014.1667
$ws.Range("M107").Value = 1698
$ws.Range("N107").Value = -4854.1667

$ws.Range("H132").Value = 12037.078
$ws.Range("I132").Value = 15964.543
$ws.Range("K132").Value = 47893.629
$ws.Range("M132").Value = -45363.629

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 980.6111
$ws.Range("J5").Value = 1085.1
$ws.Range("L5").Value = 3255.3
$ws.Range("N5").Value = -3479.3

$ws.Range("H7").Value = 40
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H131").Value = 770.5612
$ws.Range("J131").Value = 770.5612
$ws.Range("L131").Value = 2311.6836
$ws.Range("N131").Value = -12391.6836

$ws.Range("H135").Value = 980.6111
$ws.Range("J135").Value = 1085.1
$ws.Range("L135").Value = 9765.9
$ws.Range("N135").Value = -14835.9

$ws.Range("H137").Value = 27781678
$ws.Range("J137").Value = 37041916
$ws.Range("L137").Value = 111125748
$ws.Range("N137").Value = -111135948

$ws.Range("H139").Value = 3025.25
$ws.Range("I139").Value = 1486.6923
$ws.Range("J139").Value = 5882.5713
$ws.Range("K139").Value = 4460.0769
$ws.Range("L139").Value = 17647.7139
$ws.Range("M139").Value = 679.9231
$ws.Range("N139").Value = -27927.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 101059.8
$ws.Range("J10").Value = 101059.8
$ws.Range("L10").Value = 101059.8
$ws.Range("N10").Value = -101397.8

$ws.Range("H43").Value = 10000
$ws.Range("J43").Value = 10000
$ws.Range("L43").Value = 10000
$ws.Range("N43").Value = -10302

$ws.Range("H80").Value = 3312.6667
$ws.Range("I80").Value = 2932.4
$ws.Range("K80").Value = 2932.4
$ws.Range("M80").Value = -1934.4

$ws.Range("H83").Value = 3312.6667
$ws.Range("I83").Value = 2932.4
$ws.Range("K83").Value = 14662
$ws.Range("M83").Value = -9670

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 2503
$ws.Range("I4").Value = 1754.5
$ws.Range("K4").Value = 1754.5
$ws.Range("M4").Value = -1641.5

$ws.Range("H14").Value = 2297.1428
$ws.Range("J14").Value = 2297.1428
$ws.Range("L14").Value = 2297.1428
$ws.Range("N14").Value = -2641.1428

$ws.Range("H28").Value = 2503
$ws.Range("I28").Value = 1754.5
$ws.Range("K28").Value = 1754.5
$ws.Range("M28").Value = -1522.5

$ws.Range("H29").Value = 50007500

$ws.Range("H35").Value = 20999.75
$ws.Range("I35").Value = 17999.666
$ws.Range("J35").Value = 30000
$ws.Range("K35").Value = 17999.666
$ws.Range("L35").Value = 30000
$ws.Range("M35").Value = -17663.666
$ws.Range("N35").Value = -30672

$ws.Range("H37").Value = 2503
$ws.Range("I37").Value = 1754.5
$ws.Range("K37").Value = 1754.5
$ws.Range("M37").Value = -1647.5

$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4540

$ws.Range("H93").Value = 2314.8948
$ws.Range("I93").Value = 2162.7856
$ws.Range("J93").Value = 2740.8
$ws.Range("K93").Value = 2162.7856
$ws.Range("L93").Value = 2740.8
$ws.Range("M93").Value = -914.7856000000002
$ws.Range("N93").Value = -5236.8

$ws.Range("H94").Value = 40000
$ws.Range("J94").Value = 40000
$ws.Range("L94").Value = 40000
$ws.Range("N94").Value = -41352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 126
$ws.Range("I100").Value = 115.2
$ws.Range("J100").Value = 180
$ws.Range("K100").Value = 230.4
$ws.Range("L100").Value = 360
$ws.Range("M100").Value = 310.6
$ws.Range("N100").Value = -1442

$ws.Range("H136").Value = 47621732
$ws.Range("I136").Value = 76925500
$ws.Range("J136").Value = 3101.25
$ws.Range("K136").Value = 230776500
$ws.Range("L136").Value = 9303.75
$ws.Range("M136").Value = -230773950
$ws.Range("N136").Value = -14403.75
